# Switch R19 to 1210, 62R resistor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number-format/fill style used by other Mouser-part cells (e.g. H33)
# onto H35 before changing its value, so it matches the "General" formatted
# green cells used elsewhere in column H.
$ws.Range("H33").Copy()
$ws.Range("H35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 35 corresponds to reference R19.
# Order of the .Value assignments below matters for the order new entries
# land in the shared-strings table, so set them in the same order the
# original author appears to have touched the cells.
$ws.Range("H35").Value = "603-RC1210FR-0762RL "
$ws.Range("L35").Value = "Resistor 62R, 500mW, 75V, +/-1%, SMD 1210"
$ws.Range("C35").Value = "RC1210FR-0762RL"
$ws.Range("K35").Value = "Resistor_SMD:R_1210_3225Metric_Pad1.42x2.65mm_HandSolder"
$ws.Range("F35").Value = 0.122
$ws.Range("I35").ClearContents()

# Restore the selection to reflect where the edit was made.
[void]$ws.Range("O35").Select()
